# Refresh the cached "datetimeFigureOut" date field text (Date Placeholder)
# from 22/4/2021 to 25/4/2021 across the slide master and every slide layout,
# exactly as PowerPoint does for an auto-updating date footer/placeholder
# whose value is re-cached on save.

$p = $ppt.ActivePresentation
$newDate = "25/4/2021"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# Every slide layout belonging to the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}
